$wb = $excel.ActiveWorkbook

$handbackUrl = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/973e972c717b81da8d5a364dbd4accadccf9d582/e2e/379cb08c-10a4-4b40-a26e-61aa8615c179.md'
$handbackName = '379cb08c-10a4-4b40-a26e-61aa8615c179.md'
$errorDetail = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61ce5ac9f0dafa2fc6a5e1ec9c39b718128329e7/e2e/379cb08c-10a4-4b40-a26e-61aa8615c179.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/973e972c717b81da8d5a364dbd4accadccf9d582/e2e/379cb08c-10a4-4b40-a26e-61aa8615c179.md.'
$hyperlinkBlue = 15570276

# ---------------------------------------------------------------------------
# zh-cn sheet: newly generated handback report for row 6
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("I1:J1").ColumnWidth = 39.15
$ws.Range("P1").ColumnWidth = 39.15

$ws.Hyperlinks.Add($ws.Range("I6"), $handbackUrl, "", "", $handbackName)
$ws.Range("I6").Font.Underline = $true
$ws.Range("I6").Font.Color = $hyperlinkBlue

$ws.Range("J6").Value = "379cb08c-10a4-4b40-a26e-61aa8615c179.87389215bb379011025f3a9828b47094682982fc.zh-cn.xlf"

$ws.Range("K6").Value = "2016-11-14 17:39:06"
$ws.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Range("P6").Value = $errorDetail

# ---------------------------------------------------------------------------
# de-de sheet: refresh earlier report timestamps (shared string reindex) and
# generate the new handback report for row 6
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Range("I1:J1").ColumnWidth = 39.15
$ws2.Range("P1").ColumnWidth = 39.15

$ws2.Range("G2").Value = "10148e50-d50c-4834-8c42-deb3eff122ef.fde47cb0987ea16b58330b7ba4ebe365cff9f21c.de-de.xlf"
$ws2.Range("G3").Value = "2eae5ab6-36be-4338-8904-23370245bf84.3aa9fd1bef469987b993cc6cc4ec260f17572723.de-de.xlf"
$ws2.Range("G4").Value = "5e69d221-bb64-4fb3-8393-616150bcf521.b1e1cd21933f9cf06ae0f8fe9c62348cded2bf2a.de-de.xlf"
$ws2.Range("G5").Value = "0bf97cce-5374-4340-8678-921df3f3f590.061457546e8265f5d982ce94b846b747eaef5b02.de-de.xlf"
$ws2.Range("G6").Value = "379cb08c-10a4-4b40-a26e-61aa8615c179.87389215bb379011025f3a9828b47094682982fc.de-de.xlf"

$ws2.Hyperlinks.Add($ws2.Range("I6"), $handbackUrl, "", "", $handbackName)
$ws2.Range("I6").Font.Underline = $true
$ws2.Range("I6").Font.Color = $hyperlinkBlue

$ws2.Range("J6").Value = "379cb08c-10a4-4b40-a26e-61aa8615c179.87389215bb379011025f3a9828b47094682982fc.de-de.xlf"

$ws2.Range("K6").Value = "2016-11-14 17:39:27"
$ws2.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("P6").Value = $errorDetail
